$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price (D) cells so numeric-looking strings
# (e.g. "0.0000179", "0.810", "64.150.98") are preserved verbatim
# instead of being coerced to Excel numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.262.37'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.507.78'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.94'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.29'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.506.67'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.13'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.376'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.106.28'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.37'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000179'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.118'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.508.22'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.285.32'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.76'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.86'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '383.28'
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.568'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.650.92'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.91'
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.73'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000116'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.58'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.57'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.31'
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.22'
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.522.99'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.146'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '23.55'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.29'
$ws.Range('E38').Value = '  +1.58%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('E39').Value = '  -3.64%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.85'
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '162.62'
$ws.Range('E41').Value = '  -5.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0781'
$ws.Range('E42').Value = '  -3.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.810'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.36'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.81'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.22'
$ws.Range('E47').Value = '  -1.86%  '
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.479.28'
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('E51').Value = '  -1.64%  '
